# edit.ps1 - apply proofing/text edits described in the commit "updated proposal
# and presentation" to cert-mbeddr.pptx (slides "Method" and "Developers and Users").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 4 ("Method") - Content Placeholder 2
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shpMethod = $s4.Shapes.Item(2)
$trMethod = $shpMethod.TextFrame.TextRange

# WP2 bullet: "requirement" -> "requirements"
$paraWP2 = $trMethod.Paragraphs(3)
$oldWP2 = ": Traceability reports (requirement down to C code)"
$newWP2 = ": Traceability reports (requirements down to C code)"
$offWP2 = $paraWP2.Text.IndexOf($oldWP2)
$runWP2 = $trMethod.Characters($paraWP2.Start + $offWP2, $oldWP2.Length)
$runWP2.Text = $newWP2

# WP4 bullet: "(qualification?)" -> "(qualification??)"
$paraWP4 = $trMethod.Paragraphs(5)
$oldWP4 = ": Integrate all these technical arguments into a logical certification (qualification?) argument of the C code generated by "
$newWP4 = ": Integrate all these technical arguments into a logical certification (qualification??) argument of the C code generated by "
$offWP4 = $paraWP4.Text.IndexOf($oldWP4)
$runWP4 = $trMethod.Characters($paraWP4.Start + $offWP4, $oldWP4.Length)
$runWP4.Text = $newWP4

# ---------------------------------------------------------------------------
# Slide 5 ("Developers and Users") - left "Content Placeholder 2" (developers)
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shpDevs = $s5.Shapes.Item(2)
$trDevs = $shpDevs.TextFrame.TextRange

# "Itemis (with financing?)" -> "Itemis (with financing??)"
$paraItemis = $trDevs.Paragraphs(3)
$oldFin = " (with financing?)"
$newFin = " (with financing??)"
$offFin = $paraItemis.Text.IndexOf($oldFin)
$runFin = $trDevs.Characters($paraItemis.Start + $offFin, $oldFin.Length)
$runFin.Text = $newFin

# "the mbeddr team" -> italicize "mbeddr" and the following space, keep "team" upright
$paraTeam = $trDevs.Paragraphs(5)
$teamText = $paraTeam.Text
$mbIdx = $teamText.IndexOf("mbeddr")
$runMbeddr = $trDevs.Characters($paraTeam.Start + $mbIdx, 6)
$runMbeddr.Font.Italic = $true
$runSpace = $trDevs.Characters($paraTeam.Start + $mbIdx + 6, 1)
$runSpace.Font.Italic = $true

# ---------------------------------------------------------------------------
# Slide 5 ("Developers and Users") - right "Content Placeholder 2" (users)
# ---------------------------------------------------------------------------
$shpUsers = $s5.Shapes.Item(5)
$trUsers = $shpUsers.TextFrame.TextRange

# "OHB? (Aerospace)" -> "OHB (Aerospace) " + "(with financing??)" (two runs)
$paraOHB = $trUsers.Paragraphs(2)
$paraOHB.Text = "OHB (Aerospace) (with financing??)"

$ohbFull = $paraOHB.Text
$splitAt = $ohbFull.IndexOf("(with financing??)")
$runOHB1 = $trUsers.Characters($paraOHB.Start, $splitAt)
$runOHB1.Text = "OHB (Aerospace) "
$runOHB2 = $trUsers.Characters($paraOHB.Start + $splitAt, "(with financing??)".Length)
$runOHB2.Text = "(with financing??)"
